# Rebuild the shared-string cell content as pretty-printed JSON (with actual
# curly quotes, matching \u2018 / \u2019 escapes once re-serialised) and drop
# the now-unused header row (which carried the bold/bordered/centered style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuestionsText = @"
questions = [
    {
        "title": "The layouts in React Native are handled using Flexbox. To render elements side by side horizontally, which of the following properties would you use?",
        "ques_type": 2,
        "options": [
            "flex: \u2018row\u2019",
            "flexDirection: \u2018row\u2019",
            "flexbox: \u2018row\u2019",
            "justifyContent: true"
        ],
        "score": "flexDirection: \u2018row\u2019"
    },
    {
        "title": "To reuse the component logic in React, you can use a higher-order component. What is a higher-order component?",
        "ques_type": 2,
        "options": [
            "A function that takes a component and overrides it.",
            "A function that takes a component and extends it.",
            "A function that takes a component and returns another component.",
            "A function that takes a component and overrides its render function."
        ],
        "score": "A function that takes a component and returns another component."
    },
    {
        "title": "Which of the following must you use to sign iOS apps before publishing?",
        "ques_type": 2,
        "options": [
            "Distribution profile, distribution certificate, and developer signature",
            "Developer profile and distribution certificate",
            "Distribution profile and distribution certificate",
            "Developer profile and developer certificate"
        ],
        "score": "Distribution profile and distribution certificate"
    },
    {
        "title": "Which of the following tools would you use to debug the performance of a React Native Android app?",
        "ques_type": 2,
        "options": [
            "systrace",
            "perfmon",
            "lighthouse",
            "chrome dev tools"
        ],
        "score": "systrace"
    }
]
"@

# Row 1 held the placeholder "0" value with the bold+border+centered style;
# removing it shifts row 2 (the real payload) up into row 1.
$ws.Rows("1:1").Delete()

# Row 1 / A1 now holds the (unstyled) shared string that used to live in A2.
# Replace its text with the pretty-printed JSON form.
$ws.Range("A1").Value2 = $newQuestionsText
